$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 2: "Create a pdf using MatLab..." -> "Create a pdf/html using MatLab..."
#   Split into three runs: "Create a pdf" | "/html" | " using "
#   and relocate the (hidden) _GoBack bookmark to sit between "/html" and " using ".
# ---------------------------------------------------------------------------

$full = $d.Content.Text
$needle2 = "Create a pdf using "
$idx2 = $full.IndexOf($needle2)
$posAfterPdf = $idx2 + ("Create a pdf").Length

# Insert the new "/html" text right after "Create a pdf".
$insRng = $d.Range($posAfterPdf, $posAfterPdf)
$insRng.InsertAfter("/html")

$posAfterHtml = $posAfterPdf + ("/html").Length

# Force a run split between "Create a pdf" and "/html" with a throwaway bookmark,
# then remove it once the split has taken effect.
$d.Bookmarks.Add("ZZZTEMPSPLITZZZ", $d.Range($posAfterPdf, $posAfterPdf))

# Move the document's _GoBack bookmark to the new edit point (right after "/html").
$d.Bookmarks.Add("_GoBack", $d.Range($posAfterHtml, $posAfterHtml))

$d.Bookmarks.Item("ZZZTEMPSPLITZZZ").Delete()

# ---------------------------------------------------------------------------
# Change 1: "Check whether the assumptions of the linear model ar" + "e being
# satisfied (make a scatterplot with a regression line)." -> merge into a
# single run (the old _GoBack bookmark that used to separate them is already
# gone, since it was relocated above).
# ---------------------------------------------------------------------------

$old1 = "Check whether the assumptions of the linear model ar" + "e being satisfied (make a scatterplot with a regression line)."
$new1 = "Check whether the assumptions of the linear model are being satisfied (make a scatterplot with a regression line)."
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

Write-Host "change1 (merge run):" $found1
Write-Host "change2 (pdf/html insert) done at index" $posAfterPdf
